$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = -0.252
$ws.Range("G2").Value = -0.02222222222222222
$ws.Range("H2").Value = -0.02222222222222222
$ws.Range("I2").Value = -0.7886178861788619
$ws.Range("J2").Value = -0.7886178861788619
$ws.Range("K2").Value = 31.1
$ws.Range("L2").Value = 8.428184281842819
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 0.016
$ws.Range("V2").Value = 0.002043422733077906
$ws.Range("W2").Value = -0.4307479224376731
$ws.Range("X2").Value = 0.3288721801406141
$ws.Range("Y2").Value = -0.7596201025782873
$ws.Range("Z2").Value = 2.179562906083871
$ws.Range("AA2").Value = -1.718842291789719
$ws.Range("AB2").Value = 0.08457913042766638
$ws.Range("AC2").Value = -1.803421422217386
$ws.Range("AD2").Value = 68.40000000000001
$ws.Range("AF2").Value = 68.40000000000001
$ws.Range("AG2").Value = 68.384
$ws.Range("AH2").Value = 0.8972845336481701
$ws.Range("AI2").Value = 2.496350364963503
$ws.Range("AJ2").Value = 0.8972629700579946
$ws.Range("AK2").Value = 2.497224656733859
$ws.Range("AL2").Value = 2.09
$ws.Range("AM2").Value = 1.997
$ws.Range("AN2").Value = -834.1463414634146
$ws.Range("AO2").Value = -1.392344497607656
$ws.Range("AP2").Value = -833.9512195121951
$ws.Range("AQ2").Value = -1.457185778668002

# Row 3 updates
$ws.Range("D3").Value = -0.252
$ws.Range("G3").Value = -0.02222222222222222
$ws.Range("H3").Value = -0.02222222222222222
$ws.Range("I3").Value = -0.7886178861788619
$ws.Range("J3").Value = -0.7886178861788619
$ws.Range("K3").Value = 31.1
$ws.Range("L3").Value = 8.428184281842819
$ws.Range("O3").Value = -0
$ws.Range("R3").Value = -0
$ws.Range("U3").Value = 0.016
$ws.Range("V3").Value = 0.002043422733077906
$ws.Range("W3").Value = -0.4307479224376731
$ws.Range("X3").Value = 0.3288721801406141
$ws.Range("Y3").Value = -0.7596201025782873
$ws.Range("Z3").Value = 2.179562906083871
$ws.Range("AA3").Value = -1.718842291789719
$ws.Range("AB3").Value = 0.08457913042766638
$ws.Range("AC3").Value = -1.803421422217386
$ws.Range("AD3").Value = 68.40000000000001
$ws.Range("AF3").Value = 68.40000000000001
$ws.Range("AG3").Value = 68.384
$ws.Range("AH3").Value = 0.8972845336481701
$ws.Range("AI3").Value = 2.496350364963503
$ws.Range("AJ3").Value = 0.8972629700579946
$ws.Range("AK3").Value = 2.497224656733859
$ws.Range("AL3").Value = 2.09
$ws.Range("AM3").Value = 1.997
$ws.Range("AN3").Value = -834.1463414634146
$ws.Range("AO3").Value = -1.392344497607656
$ws.Range("AP3").Value = -833.9512195121951
$ws.Range("AQ3").Value = -1.457185778668002
